$wb = $excel.ActiveWorkbook

# Add a new worksheet "ODI Batting Extra" placed after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# ---- Header row (bold, centered, top-aligned, thin border - matches the
#      other sheets' header style) ----
$headerRange = $newSheet.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# Columns A, C, D, E, F hold text-like data (numbers-as-text, percents as
# text, etc.) so force a text number format before assigning so Excel does
# not reinterpret the values. Column B (BATTING_POSITION) is a real number.
$textCells = @("A2","C2","D2","E2","F2","A3","C3","D3","E3","F3","A4","C4","D4","E4","F4")
foreach ($addr in $textCells) {
    $newSheet.Range($addr).NumberFormat = "@"
}

# Row 2
$newSheet.Range("A2").Value = "4698"
$newSheet.Range("B2").Value = 4
$newSheet.Range("C2").Value = "0"
$newSheet.Range("D2").Value = "0"
$newSheet.Range("E2").Value = ""
$newSheet.Range("F2").Value = "NO"

# Row 3
$newSheet.Range("A3").Value = "4699"
$newSheet.Range("B3").Value = 4
$newSheet.Range("C3").Value = "7"
$newSheet.Range("D3").Value = "4"
$newSheet.Range("E3").Value = "23.39%"
$newSheet.Range("F3").Value = "NO"

# Row 4
$newSheet.Range("A4").Value = "4700"
$newSheet.Range("B4").Value = 4
$newSheet.Range("C4").Value = "1"
$newSheet.Range("D4").Value = "0"
$newSheet.Range("E4").Value = "1.73%"
$newSheet.Range("F4").Value = "NO"

# Restore the originally-active sheet/tab selection
$wb.Worksheets.Item(1).Activate()
